$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I2').Value = '_CROSSREF'
$ws.Range('E3').Value = '[Wei-jie%Guan%NULL%14,                  Zheng-yi%Ni%NULL%5,                  Yu%Hu%NULL%4,                  Wen-hua%Liang%NULL%4,                  Chun-quan%Ou%NULL%4,                  Jian-xing%He%NULL%4,                  Lei%Liu%NULL%8,                  Hong%Shan%NULL%6,                  Chun-liang%Lei%NULL%4,                  David S.C.%Hui%NULL%3,                  Bin%Du%NULL%3,                  Lan-juan%Li%NULL%3,                  Guang%Zeng%NULL%3,                  Kwok-Yung%Yuen%NULL%6,                  Ru-chong%Chen%NULL%4,                  Chun-li%Tang%NULL%4,                  Tao%Wang%NULL%7,                  Ping-yan%Chen%NULL%4,                  Jie%Xiang%NULL%6,                  Shi-yue%Li%NULL%4,                  Jin-lin%Wang%NULL%3,                  Zi-jing%Liang%NULL%3,                  Yi-xiang%Peng%NULL%4,                  Li%Wei%NULL%4,                  Yong%Liu%NULL%4,                  Ya-hua%Hu%NULL%4,                  Peng%Peng%NULL%7,                  Jian-ming%Wang%NULL%4,                  Ji-yang%Liu%NULL%4,                  Zhong%Chen%NULL%4,                  Gang%Li%NULL%6,                  Zhi-jian%Zheng%NULL%4,                  Shao-qin%Qiu%NULL%4,                  Jie%Luo%NULL%4,                  Chang-jiang%Ye%NULL%4,                  Shao-yong%Zhu%NULL%4,                  Nan-shan%Zhong%NULL%5]'
$ws.Range('I3').Value = '_PMC_CROSSREF'
$ws.Range('E4').Value = '[Ning%Tang%NULL%0,                  Dengju%Li%NULL%2,                  Xiong%Wang%NULL%2,                  Ziyong%Sun%NULL%4]'
$ws.Range('I4').Value = '_PMC_elsevier_CROSSREF'
$ws.Range('E5').Value = '[Akdis%Cezmi A%coreGivesNoEmail%0,                Cao%Yi-yuan%coreGivesNoEmail%3,                Dong%Xiang%coreGivesNoEmail%3,                Gao%Ya-dong%coreGivesNoEmail%3,                Yan%You-qin%coreGivesNoEmail%3,                Yang%Yi-bin%coreGivesNoEmail%3,                Yuan%Ya-dong%coreGivesNoEmail%3,                Zhang%Jin-jin%coreGivesNoEmail%3]'
$ws.Range('I5').Value = '_CROSSREF'
$ws.Range('E6').Value = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,                Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%4,                Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%4,                Subert-Salas%Lizandra%coreGivesNoEmail%4]'
$ws.Range('I6').Value = '_CROSSREF'
$ws.Range('E7').Value = '[Huijun%Chen%NULL%3,                  Juanjuan%Guo%NULL%3,                  Chen%Wang%NULL%3,                  Fan%Luo%NULL%2,                  Xuechen%Yu%NULL%3,                  Wei%Zhang%NULL%6,                  Jiafu%Li%NULL%3,                  Dongchi%Zhao%NULL%2,                  Dan%Xu%NULL%3,                  Qing%Gong%NULL%2,                  Jing%Liao%NULL%2,                  Huixia%Yang%yanghuixia@bjmu.edu.cn%3,                  Wei%Hou%houwei@whu.edu.cn%2,                  Yuanzhen%Zhang%zhangyuanzhen@whu.edu.cn%3]'
$ws.Range('I7').Value = '_PMC_elsevier_CROSSREF'
$ws.Range('E8').Value = '[Chaolin%Huang%NULL%0,                  Yeming%Wang%NULL%6,                  Xingwang%Li%NULL%7,                  Lili%Ren%NULL%4,                  Jianping%Zhao%NULL%9,                  Yi%Hu%NULL%7,                  Li%Zhang%NULL%7,                  Guohui%Fan%NULL%6,                  Jiuyang%Xu%NULL%6,                  Xiaoying%Gu%NULL%6,                  Zhenshun%Cheng%NULL%4,                  Ting%Yu%NULL%13,                  Jiaan%Xia%NULL%4,                  Yuan%Wei%NULL%10,                  Wenjuan%Wu%NULL%4,                  Xuelei%Xie%NULL%4,                  Wen%Yin%NULL%6,                  Hui%Li%NULL%7,                  Min%Liu%NULL%4,                  Yan%Xiao%NULL%5,                  Hong%Gao%NULL%6,                  Li%Guo%NULL%5,                  Jungang%Xie%NULL%5,                  Guangfa%Wang%NULL%4,                  Rongmeng%Jiang%NULL%4,                  Zhancheng%Gao%NULL%5,                  Qi%Jin%NULL%5,                  Jianwei%Wang%wangjw28@163.com%0,                  Bin%Cao%caobin_ben@163.com%7]'
$ws.Range('I8').Value = '_PMC_elsevier_CROSSREF'
$ws.Range('E9').Value = '[Xiao-Wei%Xu%NULL%0,                  Xiao-Xin%Wu%NULL%4,                  Xian-Gao%Jiang%NULL%4,                  Kai-Jin%Xu%NULL%4,                  Ling-Jun%Ying%NULL%4,                  Chun-Lian%Ma%NULL%4,                  Shi-Bo%Li%NULL%4,                  Hua-Ying%Wang%NULL%4,                  Sheng%Zhang%NULL%4,                  Hai-Nv%Gao%NULL%4,                  Ji-Fang%Sheng%NULL%4,                  Hong-Liu%Cai%NULL%4,                  Yun-Qing%Qiu%NULL%4,                  Lan-Juan%Li%NULL%4]'
$ws.Range('I9').Value = '_PMC_CROSSREF'
$ws.Range('E10').Value = '[Chaomin%Wu%xref no email%0,            Xiaoyan%Chen%xref no email%8,            Yanping%Cai%xref no email%8,            Jia\u2019an%Xia%xref no email%8,            Xing%Zhou%xref no email%8,            Sha%Xu%xref no email%8,            Hanping%Huang%xref no email%8,            Li%Zhang%xref no email%8,            Xia%Zhou%xref no email%8,            Chunling%Du%xref no email%8,            Yuye%Zhang%xref no email%8,            Juan%Song%xref no email%9,            Sijiao%Wang%xref no email%8,            Yencheng%Chao%xref no email%8,            Zeyong%Yang%xref no email%8,            Jie%Xu%xref no email%8,            Xin%Zhou%xref no email%8,            Dechang%Chen%xref no email%8,            Weining%Xiong%xref no email%8,            Lei%Xu%xref no email%8,            Feng%Zhou%xref no email%8,            Jinjun%Jiang%xref no email%8,            Chunxue%Bai%xref no email%8,            Junhua%Zheng%xref no email%8,            Yuanlin%Song%xref no email%9]'
$ws.Range('I10').Value = '_CROSSREF'
$ws.Range('E11').Value = '[Wanbo%Zhu%NULL%0,                  Kai%Xie%NULL%3,                  Hui%Lu%NULL%2,                  Lei%Xu%bayinhexl@126.com%2,                  Shusheng%Zhou%zhouss108@163.com%2,                  Shiyuan%Fang%fangshiyuan2008@126.com%2]'
$ws.Range('I11').Value = '_PMC_CROSSREF'
$ws.Range('E12').Value = '[Dahai%Zhao%NULL%0,                  Feifei%Yao%NULL%1,                  Lijie%Wang%NULL%1,                  Ling%Zheng%NULL%1,                  Yongjun%Gao%NULL%1,                  Jun%Ye%NULL%1,                  Feng%Guo%NULL%1,                  Hui%Zhao%NULL%1,                  Rongbao%Gao%rongbaogao@163.com%1]'
$ws.Range('I12').Value = '_PMC_CROSSREF'
$ws.Range('E13').Value = '[Zhiliang%Hu%NULL%0,                  Ci%Song%NULL%1,                  Chuanjun%Xu%NULL%1,                  Guangfu%Jin%NULL%1,                  Yaling%Chen%NULL%1,                  Xin%Xu%NULL%1,                  Hongxia%Ma%NULL%1,                  Wei%Chen%NULL%3,                  Yuan%Lin%NULL%1,                  Yishan%Zheng%NULL%1,                  Jianming%Wang%NULL%1,                  Zhibin%Hu%zhibin_hu@njmu.edu.cn%1,                  Yongxiang%Yi%ian0126@126.com%1,                  Hongbing%Shen%NULL%1]'
$ws.Range('I13').Value = '_PMC_Springer_CROSSREF'
$ws.Range('E14').Value = '[Tao%Chen%NULL%12,                  Di%Wu%NULL%2,                  Huilong%Chen%NULL%2,                  Weiming%Yan%NULL%2,                  Danlei%Yang%NULL%2,                  Guang%Chen%NULL%2,                  Ke%Ma%NULL%2,                  Dong%Xu%NULL%2,                  Haijing%Yu%NULL%2,                  Hongwu%Wang%NULL%2,                  Tao%Wang%NULL%0,                  Wei%Guo%NULL%2,                  Jia%Chen%NULL%2,                  Chen%Ding%NULL%2,                  Xiaoping%Zhang%NULL%2,                  Jiaquan%Huang%NULL%2,                  Meifang%Han%NULL%2,                  Shusheng%Li%NULL%2,                  Xiaoping%Luo%NULL%2,                  Jianping%Zhao%NULL%0,                  Qin%Ning%NULL%2]'
$ws.Range('I14').Value = '_PMC_CROSSREF'
$ws.Range('E15').Value = '[Kui%Liu%NULL%0,                  Yuan-Yuan%Fang%NULL%4,                  Yan%Deng%NULL%3,                  Wei%Liu%NULL%5,                  Mei-Fang%Wang%NULL%3,                  Jing-Ping%Ma%NULL%3,                  Wei%Xiao%NULL%3,                  Ying-Nan%Wang%NULL%3,                  Min-Hua%Zhong%NULL%3,                  Cheng-Hong%Li%NULL%3,                  Guang-Cai%Li%NULL%3,                  Hui-Guo%Liu%NULL%4,                  Xiu-Yuan%Hao%NULL%6,                  Pei-Fang%Wei%NULL%8]'
$ws.Range('I15').Value = '_PMC_CROSSREF'
$ws.Range('E16').Value = '[Feng%Pan%NULL%4,                  Tianhe%Ye%NULL%3,                  Peng%Sun%NULL%3,                  Shan%Gui%NULL%3,                  Bo%Liang%NULL%3,                  Lingli%Li%NULL%3,                  Dandan%Zheng%NULL%3,                  Jiazheng%Wang%NULL%3,                  Richard L.%Hesketh%NULL%3,                  Lian%Yang%yanglian@hust.edu.cn%3,                  Chuansheng%Zheng%NULL%5]'
$ws.Range('I16').Value = '_PMC_CROSSREF'
$ws.Range('E18').Value = '[Xiaobo%Yang%NULL%0,                  Yuan%Yu%NULL%2,                  Jiqian%Xu%NULL%1,                  Huaqing%Shu%NULL%1,                  Jia''an%Xia%NULL%4,                  Hong%Liu%NULL%3,                  Yongran%Wu%NULL%2,                  Lu%Zhang%NULL%1,                  Zhui%Yu%NULL%1,                  Minghao%Fang%NULL%1,                  Ting%Yu%NULL%0,                  Yaxin%Wang%NULL%2,                  Shangwen%Pan%NULL%1,                  Xiaojing%Zou%NULL%2,                  Shiying%Yuan%NULL%2,                  You%Shang%NULL%2]'
$ws.Range('I18').Value = '_PMC_elsevier_CROSSREF'
$ws.Range('E19').Value = '[Jiaojiao%Chu%NULL%3,                  Nan%Yang%NULL%2,                  Yanqiu%Wei%NULL%2,                  Huihui%Yue%NULL%2,                  Fengqin%Zhang%NULL%2,                  Jianping%Zhao%NULL%0,                  Li%He%NULL%2,                  Gaohong%Sheng%NULL%2,                  Peng%Chen%NULL%2,                  Gang%Li%NULL%0,                  Sisi%Wu%NULL%2,                  Bo%Zhang%NULL%2,                  Shu%Zhang%szhang@tjh.tjmu.edu.cn%2,                  Congyi%Wang%wangcy@tjh.tjmu.edu.cn%2,                  Xiaoping%Miao%miaoxp@hust.edu.cn%2,                  Juan%Li%lijuan@tjh.tjmu.edu.cn%2,                  Wenhua%Liu%liuwh_2013@126.com%2,                  Huilan%Zhang%huilanz_76@163.com%2]'
$ws.Range('I19').Value = '_PMC_CROSSREF'
$ws.Range('E20').Value = '[ J.%Liu%null%1,                   L.% Ouyang%null%1,                   P.% Guo%null%1,                   H. s.% Wu%null%1,                   P.% Fu%null%1,                   Y. l.% Chen%null%1,                   D.% Yang%null%1,                   X. y.% Han%null%1,                   Y. k.% Cao%null%1,                   O.% Alwalid%null%1,                   J.% Tao%null%1,                   S. y.% Peng%null%1,                   H. s.% Shi%null%1,                   F.% Yang%null%1,                   C. s. % Zheng%null%1]'
$ws.Range('I20').Value = '_MedBiorxiv_CROSSREF'
$ws.Range('E21').Value = '[ R.%Liu%null%1,                   X.% Ming%null%1,                   O.% Xu%null%1,                   J.% Zhou%null%1,                   H.% Peng%null%1,                   N.% Xiang%null%1,                   J.% Zhang%null%1,                   H. % Zhu%null%1]'
$ws.Range('I21').Value = '_MedBiorxiv_CROSSREF'
$ws.Range('E22').Value = '[Jiaxun%Hu%xref no email%1,            You%Lv%xref no email%1,            Ying%Xu%xref no email%1,            Yushan%Miao%xref no email%1,            Wei%Wang%xref no email%5,            Wenqing%Le%xref no email%1,            Hao%Tang%xref no email%1]'
$ws.Range('I22').Value = '_CROSSREF'
$ws.Range('E23').Value = '[Temet M.%McMichael%NULL%2,                  Shauna%Clark%NULL%2,                  Sargis%Pogosjans%NULL%2,                  Meagan%Kay%NULL%2,                  James%Lewis%NULL%3,                  Atar%Baer%NULL%2,                  Vance%Kawakami%NULL%2,                  Margaret D.%Lukoff%NULL%2,                  Jessica%Ferro%NULL%2,                  Claire%Brostrom-Smith%NULL%2,                  Francis X.%Riedo%NULL%2,                  Denny%Russell%NULL%2,                  Brian%Hiatt%NULL%2,                  Patricia%Montgomery%NULL%2,                  Agam K.%Rao%NULL%2,                  Dustin W.%Currie%NULL%2,                  Eric J.%Chow%NULL%2,                  Farrell%Tobolowsky%NULL%2,                  Ana C.%Bardossy%NULL%2,                  Lisa P.%Oakley%NULL%2,                  Jesica R.%Jacobs%NULL%2,                  Noah G.%Schwartz%NULL%2,                  Nimalie%Stone%NULL%1,                  Sujan C.%Reddy%NULL%4,                  John A.%Jernigan%NULL%3,                  Margaret A.%Honein%NULL%2,                  Thomas A.%Clark%NULL%2,                  Jeffrey S.%Duchin%NULL%2,                  NULL%NULL%NULL%22,                  NULL%NULL%NULL%0,                  Meaghan S.%Fagalde%NULL%1,                  Jennifer L.%Lenahan%NULL%1,                  Emily B.%Maier%NULL%1,                  Kaitlyn J.%Sykes%NULL%1,                  Grace%Hatt%NULL%1,                  Holly%Whitney%NULL%1,                  Melinda%Huntington-Frazier%NULL%1,                  Elysia%Gonzales%NULL%1,                  Laura A.%Mummert%NULL%1,                  Hal Garcia%Smith%NULL%1,                  Steve%Stearns%NULL%1,                  Eileen%Benoliel%NULL%1,                  Shelly%McKeirnan%NULL%1,                  Jennifer L.%Morgan%NULL%1,                  Daniel%Smith%NULL%1,                  Michaela%Hope%NULL%1,                  Noel%Hatley%NULL%1,                  Leslie M.%Barnard%NULL%1,                  Leilani%Schwarcz%NULL%1,                  Nicole%Yarid%NULL%1,                  Eric%Yim%NULL%1,                  Sandra%Kreider%NULL%1,                  Dawn%Barr%NULL%1,                  Nancy%Wilde%NULL%1,                  Courtney%Dorman%NULL%1,                  Airin%Lam%NULL%1,                  Jeanette%Harris%NULL%1,                  Hollianne%Bruce%NULL%2,                  Christopher%Spitters%NULL%2,                  Snohomish Health%District%NULL%1,                  Rachael%Zacks%NULL%1,                  Jonathan%Dyal%NULL%1,                  Michael%Hughes%NULL%1,                  Christina%Carlson%NULL%1,                  Barbara%Cooper%NULL%1,                  Michelle%Banks%NULL%1,                  Heather%McLaughlin%NULL%1,                  Arun%Balajee%NULL%1,                  Christine%Olson%NULL%1,                  Suzanne%Zane%NULL%1,                  Hammad%Ali%NULL%1,                  Jessica%Healy%NULL%1,                  Kristine%Schmit%NULL%1,                  Kevin%Spicer%NULL%1,                  Zeshan%Chisty%NULL%1,                  Sukarma%Tanwar%NULL%1,                  Joanne%Taylor%NULL%1,                  Leisha%Nolen%NULL%1,                  Jeneita%Bell%NULL%1,                  Kelly%Hatfield%NULL%1,                  Melissa%Arons%NULL%1,                  Anne%Kimball%NULL%1,                  Allison%James%NULL%1,                  Mark%Methner%NULL%1,                  Joshua%Harney%NULL%1]'
